$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row cells: "_old" suffix -> "_FV2404", "_new" suffix -> "_FV2410"
$headerMap = @{
    "A1" = "Segmentname_FV2404";
    "B1" = "Segmentgruppe_FV2404";
    "C1" = "Segment_FV2404";
    "D1" = "Datenelement_FV2404";
    "E1" = "Segment ID_FV2404";
    "F1" = "Code_FV2404";
    "G1" = "Qualifier_FV2404";
    "H1" = "Beschreibung_FV2404";
    "I1" = "Bedingungsausdruck_FV2404";
    "J1" = "Bedingung_FV2404";
    "L1" = "Segmentname_FV2410";
    "M1" = "Segmentgruppe_FV2410";
    "N1" = "Segment_FV2410";
    "O1" = "Datenelement_FV2410";
    "P1" = "Segment ID_FV2410";
    "Q1" = "Code_FV2410";
    "R1" = "Qualifier_FV2410";
    "S1" = "Beschreibung_FV2410";
    "T1" = "Bedingungsausdruck_FV2410";
    "U1" = "Bedingung_FV2410";
}

foreach ($addr in $headerMap.Keys) {
    $ws.Range($addr).Value = $headerMap[$addr]
}

# 2. Freeze panes: freeze the first row
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3. Convert the data range into an Excel Table (ListObject)
$range = $ws.Range("A1:U72")
$listObject = $ws.ListObjects.Add(1, $range, $null, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = "TableStyleMedium9"
